$wb = $excel.ActiveWorkbook

# =========================================================================
# Sheet 1 - the Key/Value info box becomes a "Table 1" banner over a
# Date/Time/Place/Latitude/Longitude/Timezone/Sunrise/Sunset/Ayanamsha
# block (everything below the old header shifts up by one row, and the
# Name row is dropped while the Date gets a new value).
# =========================================================================
$ws1 = $wb.Worksheets.Item("Sheet 1")

# Row 1 becomes a single "Table 1" banner cell (will be merged below).
$ws1.Range("A1").Value = "Table 1"

# Rows 2-10: shift every label/value up by one row and refresh with the new data set.
$ws1.Range("A2").Value = "Date"
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "09/09/1989"
$ws1.Range("B2").ClearFormats() | Out-Null

$ws1.Range("A3").Value = "Time"
$ws1.Range("B3").Value = "03:04:00"

$ws1.Range("A4").Value = "Place"
$ws1.Range("B4").Value = "Puttur,Karnataka"

$ws1.Range("A5").Value = "Latitude"
$ws1.Range("B5").Value = 12.7632858

$ws1.Range("A6").Value = "Longitude"
$ws1.Range("B6").Value = 75.20184209999999

$ws1.Range("A7").Value = "Timezone"
$ws1.Range("B7").Value = "tobefilled GMT+5.5"

$ws1.Range("A8").Value = "Sunrise"
$ws1.Range("B8").Value = "tobefilled"

$ws1.Range("A9").Value = "Sunset"
$ws1.Range("B9").Value = "tobefilled"

$ws1.Range("A10").Value = "Ayanamsha"
$ws1.Range("B10").Value = "tobefilled"

# The old row 11 (Ayanamsha/tobefilled) is gone now that everything moved up one row.
$ws1.Range("A11:B11").Clear() | Out-Null

# Merge the new banner row across the table's full width (A:K, same width as Sheet 2's
# banner), then drop the bold/bordered style the old A1 header cell used to carry.
$ws1.Range("A1:K1").Merge() | Out-Null
$ws1.Range("A1").ClearFormats() | Out-Null

# =========================================================================
# Sheet 2 - planetary position table: recompute the per-planet Sign, Sign
# Lord, Nakshatra, Naksh Lord, Degree, Retro, Combust and House columns
# for the new chart (column A - the planet names - is unchanged).
# =========================================================================
$ws2 = $wb.Worksheets.Item("Sheet 2")

$ws2.Range("B3").Value = "Leo"
$ws2.Range("C3").Value = "Sun"
$ws2.Range("D3").Value = "Purva Phalguni"
$ws2.Range("F3").Value = 142.4845860305107

$ws2.Range("B4").Value = "Scorpio"
$ws2.Range("C4").Value = "Mars"
$ws2.Range("D4").Value = "Jyeshtha"
$ws2.Range("E4").Value = "Mercury"
$ws2.Range("F4").Value = 238.080245280088
$ws2.Range("J4").Value = 5

$ws2.Range("B5").Value = "Virgo"
$ws2.Range("C5").Value = "Mercury"
$ws2.Range("D5").Value = "Hasta"
$ws2.Range("E5").Value = "Moon"
$ws2.Range("F5").Value = 166.4089687091878
$ws2.Range("J5").Value = 3

$ws2.Range("B6").Value = "Libra"
$ws2.Range("C6").Value = "Venus"
$ws2.Range("D6").Value = "Chitra"
$ws2.Range("E6").Value = "Mars"
$ws2.Range("F6").Value = 182.0572724410369

$ws2.Range("B7").Value = "Leo"
$ws2.Range("C7").Value = "Sun"
$ws2.Range("D7").Value = "Uttara Phalguni"
$ws2.Range("E7").Value = "Sun"
$ws2.Range("F7").Value = 149.3837157754606
$ws2.Range("G7").Value = "Direct"
$ws2.Range("H7").Value = "Combust"
$ws2.Range("J7").Value = 2

$ws2.Range("B8").Value = "Gemini"
$ws2.Range("C8").Value = "Mercury"
$ws2.Range("D8").Value = "Ardra"
$ws2.Range("E8").Value = "Rahu"
$ws2.Range("F8").Value = 73.33655250974485
$ws2.Range("J8").Value = 12

$ws2.Range("B9").Value = "Sagittarius"
$ws2.Range("C9").Value = "Jupiter"
$ws2.Range("D9").Value = "Purva Ashadha"
$ws2.Range("F9").Value = 253.590182561376
$ws2.Range("G9").Value = "Retro"
$ws2.Range("H9").Value = "No"
$ws2.Range("J9").Value = 6

$ws2.Range("B10").Value = "Sagittarius"
$ws2.Range("C10").Value = "Jupiter"
$ws2.Range("D10").Value = "Moola"
$ws2.Range("E10").Value = "Ketu"
$ws2.Range("F10").Value = 247.622252771897
$ws2.Range("G10").Value = "Retro"
$ws2.Range("J10").Value = 6

$ws2.Range("B11").Value = "Sagittarius"
$ws2.Range("C11").Value = "Jupiter"
$ws2.Range("D11").Value = "Purva Ashadha"
$ws2.Range("E11").Value = "Venus"
$ws2.Range("F11").Value = 255.9345667638737
$ws2.Range("J11").Value = 6

$ws2.Range("B12").Value = "Libra"
$ws2.Range("C12").Value = "Venus"
$ws2.Range("D12").Value = "Swati"
$ws2.Range("E12").Value = "Rahu"
$ws2.Range("F12").Value = 199.2862020569684
$ws2.Range("G12").Value = "Direct"
$ws2.Range("J12").Value = 4

$ws2.Range("B13").Value = "Aquarius"
$ws2.Range("C13").Value = "Saturn"
$ws2.Range("D13").Value = "Dhanishta"
$ws2.Range("E13").Value = "Mars"
$ws2.Range("F13").Value = 300.7872291398033
$ws2.Range("J13").Value = 8

$ws2.Range("B14").Value = "Leo"
$ws2.Range("C14").Value = "Sun"
$ws2.Range("D14").Value = "Magha"
$ws2.Range("E14").Value = "Ketu"
$ws2.Range("F14").Value = 121.9004845901756
$ws2.Range("G14").Value = "Direct"
$ws2.Range("J14").Value = 2
